# Add four new "column" rows (bancada, proponente, prop_pk1, prop_pk2)
# to the end of the Relatorio_SISGESAC columns listing, continuing the
# existing numbering/formatting pattern used by the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$xlPasteFormats = -4122

# Copy the formatting of the last existing "index" cell (A60) down onto
# the four new index cells so they pick up the same style (s="1") as the
# rest of column A.
$ws.Range("A60").Copy() | Out-Null
$ws.Range("A61:A64").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false

$newColumns = @(
    @{ Row = 61; Index = 59; Name = "bancada" },
    @{ Row = 62; Index = 60; Name = "proponente" },
    @{ Row = 63; Index = 61; Name = "prop_pk1" },
    @{ Row = 64; Index = 62; Name = "prop_pk2" }
)

foreach ($entry in $newColumns) {
    $ws.Cells.Item($entry.Row, 1).Value = $entry.Index
    $ws.Cells.Item($entry.Row, 2).Value = $entry.Name
}
